# Case_1_75 res_bus/vm_pu.xlsx update: case with 380 kV done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02
    "C2" = 1.030810476033444
    "D2" = 1.041017884586774
    "E2" = 1.040003573737339
    "F2" = 1.049861915538041
    "I2" = 1.040192023258582
    "J2" = 1.035949664630284
    "K2" = 1.043798589620965
    "L2" = 1.042787155324714
    "M2" = 1.052617794588261
    "N2" = 1.015980747531803
    "B3" = 1.02
    "C3" = 1.031637328132824
    "D3" = 1.041557981943783
    "E3" = 1.040749909746267
    "F3" = 1.050718226019307
    "I3" = 1.040382465757945
    "J3" = 1.03641875291547
    "K3" = 1.044149939817283
    "L3" = 1.043343994322152
    "M3" = 1.053286322611705
    "N3" = 1.016136685034754
    "B4" = 1.02
    "C4" = 1.032172897272404
    "D4" = 1.041907401386499
    "E4" = 1.04123371757628
    "F4" = 1.051273301520107
    "I4" = 1.040504426056724
    "J4" = 1.03672216441785
    "K4" = 1.044376518097333
    "L4" = 1.043704518394724
    "M4" = 1.053719245464106
    "N4" = 1.016237522384496
    "B5" = 1.02
    "C5" = 1.032398178325483
    "D5" = 1.042054281048676
    "E5" = 1.041437318886273
    "F5" = 1.051506889213311
    "I5" = 1.040555393497864
    "J5" = 1.036849688643903
    "K5" = 1.044471586105348
    "L5" = 1.043856132042906
    "M5" = 1.053901326095907
    "N5" = 1.016279898482912
    "B6" = 1.02
    "C6" = 1.032436011451999
    "D6" = 1.042078941778145
    "E6" = 1.04147151661551
    "F6" = 1.051546123277378
    "I6" = 1.040563933262499
    "J6" = 1.03687109872895
    "K6" = 1.04448753752209
    "L6" = 1.043881591492683
    "M6" = 1.053931902870309
    "N6" = 1.016287012669506
    "B7" = 1.02
    "C7" = 1.032175906988774
    "D7" = 1.041909364066182
    "E7" = 1.041236437289987
    "F7" = 1.051276421813467
    "I7" = 1.040505108284474
    "J7" = 1.036723868522124
    "K7" = 1.044377789131429
    "L7" = 1.043706544071679
    "M7" = 1.053721678122688
    "N7" = 1.016238088679129
    "B8" = 1.02
    "C8" = 1.031089801674747
    "D8" = 1.041200424177386
    "E8" = 1.040255618486016
    "F8" = 1.050151104480498
    "I8" = 1.040256646332838
    "J8" = 1.036108219308125
    "K8" = 1.043917488377008
    "L8" = 1.042975296916292
    "M8" = 1.052843655262002
    "N8" = 1.016033460456379
    "B9" = 1.02
    "C9" = 1.029180156877167
    "D9" = 1.039950812252628
    "E9" = 1.038534099419445
    "F9" = 1.048175784169951
    "I9" = 1.039809148760323
    "J9" = 1.035022502771278
    "K9" = 1.043100555607477
    "L9" = 1.041688432738455
    "M9" = 1.051299147928496
    "N9" = 1.015672404084091
    "B10" = 1.02
    "C10" = 1.027909988186297
    "D10" = 1.039117607695721
    "E10" = 1.037391106577728
    "F10" = 1.046864151477467
    "I10" = 1.039504366865442
    "J10" = 1.034298184863015
    "K10" = 1.042552104281001
    "L10" = 1.040831742244107
    "M10" = 1.050271377563811
    "N10" = 1.015431406906941
    "B11" = 1.02
    "C11" = 1.027360705755592
    "D11" = 1.038756815111453
    "E11" = 1.036897312552017
    "F11" = 1.046297469087917
    "I11" = 1.039370874668122
    "J11" = 1.033984442549891
    "K11" = 1.042313728691533
    "L11" = 1.040461092804794
    "M11" = 1.049826813345625
    "N11" = 1.015326988596665
    "B12" = 1.02
    "C12" = 1.027156785936697
    "D12" = 1.038622801276201
    "E12" = 1.036714067002476
    "F12" = 1.046087169774426
    "I12" = 1.039321062332869
    "J12" = 1.033867889724063
    "K12" = 1.042225052747787
    "L12" = 1.040323464091023
    "M12" = 1.049661754253169
    "N12" = 1.015288193679734
    "B13" = 1.02
    "C13" = 1.027200522510904
    "D13" = 1.038651547639443
    "E13" = 1.036753366002681
    "F13" = 1.046132270981956
    "I13" = 1.039331757523007
    "J13" = 1.033892891357885
    "K13" = 1.042244080017423
    "L13" = 1.040352983767779
    "M13" = 1.049697156710517
    "N13" = 1.01529651573507
    "B14" = 1.02
    "C14" = 1.027343847464532
    "D14" = 1.038745737462213
    "E14" = 1.036882161907688
    "F14" = 1.046280081757862
    "I14" = 1.039366761798326
    "J14" = 1.033974808554953
    "K14" = 1.042306401407676
    "L14" = 1.040449715405458
    "M14" = 1.049813168034126
    "N14" = 1.015323781983997
    "B15" = 1.02
    "C15" = 1.027432169075732
    "D15" = 1.038803771059915
    "E15" = 1.036961540101445
    "F15" = 1.046371178323986
    "I15" = 1.03938829898375
    "J15" = 1.034025278498411
    "K15" = 1.042344782141064
    "L15" = 1.040509321229784
    "M15" = 1.049884656010561
    "N15" = 1.015340580399994
    "B16" = 1.02
    "C16" = 1.027946457327649
    "D16" = 1.039141552283806
    "E16" = 1.037423902014917
    "F16" = 1.04690178710439
    "I16" = 1.039513194364471
    "J16" = 1.034319004740564
    "K16" = 1.042567905803499
    "L16" = 1.040856347537594
    "M16" = 1.05030089181144
    "N16" = 1.015438335476431
    "B17" = 1.02
    "C17" = 1.028269247717658
    "D17" = 1.039353432405707
    "E17" = 1.037714233115608
    "F17" = 1.047234963567205
    "I17" = 1.039591131737875
    "J17" = 1.034503223449911
    "K17" = 1.04270762745329
    "L17" = 1.04107411016323
    "M17" = 1.050562111777105
    "N17" = 1.015499637583907
    "B18" = 1.02
    "C18" = 1.028457594190283
    "D18" = 1.039477017442889
    "E18" = 1.037883687182724
    "F18" = 1.047429421511425
    "I18" = 1.039636444605893
    "J18" = 1.034610664614053
    "K18" = 1.042789038512368
    "L18" = 1.041201156576296
    "M18" = 1.050714521923048
    "N18" = 1.015535387769949
    "B19" = 1.02
    "C19" = 1.028521827005802
    "D19" = 1.039519156506522
    "E19" = 1.037941485046559
    "F19" = 1.047495747236097
    "I19" = 1.039651870214953
    "J19" = 1.034647297456894
    "K19" = 1.042816782911605
    "L19" = 1.041244481021841
    "M19" = 1.05076649741923
    "N19" = 1.015547576577909
    "B20" = 1.02
    "C20" = 1.028234608287492
    "D20" = 1.039330699768165
    "E20" = 1.037683072064107
    "F20" = 1.047199204278445
    "I20" = 1.039582784960819
    "J20" = 1.034483459602094
    "K20" = 1.042692645549704
    "L20" = 1.041050743268068
    "M20" = 1.050534080684518
    "N20" = 1.01549306109592
    "B21" = 1.02
    "C21" = 1.027301638840005
    "D21" = 1.038718000874387
    "E21" = 1.036844229973645
    "F21" = 1.046236549859162
    "I21" = 1.039356460181403
    "J21" = 1.033950686385152
    "K21" = 1.042288052961278
    "L21" = 1.040421229040398
    "M21" = 1.049779003593937
    "N21" = 1.0153157530066
    "B22" = 1.02
    "C22" = 1.026715670341748
    "D22" = 1.038332777200176
    "E22" = 1.03631781003929
    "F22" = 1.045632401243597
    "I22" = 1.039212845658114
    "J22" = 1.033615624772394
    "K22" = 1.042032902870065
    "L22" = 1.040025701047819
    "M22" = 1.049304673200332
    "N22" = 1.015204218746617
    "B23" = 1.02
    "C23" = 1.027026243412359
    "D23" = 1.038536990463874
    "E23" = 1.03659678042664
    "F23" = 1.045952565941682
    "I23" = 1.039289102789325
    "J23" = 1.03379325500364
    "K23" = 1.04216823494667
    "L23" = 1.040235351601988
    "M23" = 1.049556084675539
    "N23" = 1.015263350101385
    "B24" = 1.02
    "C24" = 1.028250260144102
    "D24" = 1.039340971672071
    "E24" = 1.037697152065882
    "F24" = 1.047215361985852
    "I24" = 1.039586556963423
    "J24" = 1.034492390066768
    "K24" = 1.04269941549362
    "L24" = 1.04106130167226
    "M24" = 1.050546746589961
    "N24" = 1.015496032747268
    "B25" = 1.02
    "C25" = 1.029673336832071
    "D25" = 1.040273898429275
    "E25" = 1.038978335268377
    "F25" = 1.048685535380365
    "I25" = 1.039925978545147
    "J25" = 1.03530328127612
    "K25" = 1.043312433583729
    "L25" = 1.042020909334904
    "M25" = 1.051698112208818
    "N25" = 1.015765799389391
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
